$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header G1
$ws.Range("G1").Value = "Data"

# Add rows 13-19, copies of row 2 (Urban Hero Gold / Jimmy Choo)
for ($r = 13; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = "Urban Hero Gold "
    $ws.Cells.Item($r, 2).Value = "Jimmy Choo"
    $ws.Cells.Item($r, 3).Value = "woda perfumowana dla mężczyzn"
    $ws.Cells.Item($r, 4).Value = "50 ml "
    $ws.Cells.Item($r, 5).Value = 144.5
    $ws.Cells.Item($r, 6).Value = "https://www.notino.pl/jimmy-choo/urban-hero-gold-woda-perfumowana-dla-mezczyzn/"
}

# G column: rows 14-18 get a date value (serial 44989 = 2023-03-04), row 19 gets a text date
$ws.Cells.Item(14, 7).Value = 44989
$ws.Cells.Item(15, 7).Value = 44989
$ws.Cells.Item(16, 7).Value = 44989
$ws.Cells.Item(17, 7).Value = 44989
$ws.Cells.Item(18, 7).Value = 44989
$ws.Range("G14:G18").NumberFormat = "mm-dd-yy"

$ws.Range("G19").NumberFormat = "@"
$ws.Cells.Item(19, 7).Value = "04.03.2023"

$ws.Columns.Item(6).EntireColumn.AutoFit()
$ws.Columns.Item(7).EntireColumn.AutoFit()

$ws.Range("G24").Select()
